$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 for the machine-readable ("snake-case") field
# names, shifting the existing metadata rows (iaest-measure/sdmx mappings,
# medida/dim, xsd types) down by one.
$ws.Rows(2).Insert()

# The old sparse row (previously row 5, containing only the stray
# "mapping-ano.xlsx" value at M5) has now been pushed down to row 6 - drop it,
# since the real xsd-type row (old row 4) now lives at row 5.
$ws.Rows(6).Delete()

# Populate the new row 2 with the snake-case field identifiers, which let two
# columns be related to build SKOS hierarchies.
$ws.Range("A2").Value = "inscripcion-comarca-codigo"
$ws.Range("B2").Value = "residencia-continente-nombre"
$ws.Range("C2").Value = "residencia-area-nombre"
$ws.Range("D2").Value = "personas"
$ws.Range("E2").Value = "residencia-area-codigo"
$ws.Range("F2").Value = "inscripcion-municipio-codigo"
$ws.Range("G2").Value = "sexo"
$ws.Range("H2").Value = "inscripcion-municipio-estrato"
$ws.Range("I2").Value = "inscripcion-municipio-nombre"
$ws.Range("J2").Value = "inscripcion-provincia-nombre"
$ws.Range("K2").Value = "inscripcion-provincia-codigo"
$ws.Range("L2").Value = "inscripcion-comarca-nombre"
$ws.Range("M2").Value = "ano"
